$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.639.16"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.704.82"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'314.50"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.3977"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").Value = "'0.4063"
$ws.Range("E8").Value = "  +1.83%  "
$ws.Range("D9").Value = "'1.001"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "'1.515"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("D11").Value = "'52.94"
$ws.Range("E11").Value = "  +9.31%  "
$ws.Range("D12").Value = "'0.08834"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "'7.330"
$ws.Range("E13").Value = "  +11.07%  "
$ws.Range("D14").Value = "'23.43"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "'0.00001330"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'7.585"
$ws.Range("E16").Value = "  +5.11%  "
$ws.Range("D17").Value = "1.703.90"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "'101.06"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'0.07136"
$ws.Range("E19").Value = "  +4.63%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'6.787"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'0.9988"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "'14.20"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").Value = "24.637.84"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +8.04%  "
$ws.Range("D26").Value = "'2.311"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'22.49"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").Value = "'160.46"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'5.133"
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("D30").Value = "'133.94"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "'7.367"
$ws.Range("E31").Value = "  +26.50%  "
$ws.Range("D32").Value = "1.890.02"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "'1.089"
$ws.Range("E33").Value = "  -7.13%  "
$ws.Range("D34").Value = "'0.08726"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "'7.263"
$ws.Range("E35").Value = "  +17.93%  "
$ws.Range("D36").Value = "'11.16"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("D37").Value = "'1.960"
$ws.Range("E37").Value = "  +5.92%  "
$ws.Range("D38").Value = "'0.2729"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").Value = "'14.87"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").Value = "'0.02778"
$ws.Range("E40").Value = "  +9.27%  "
$ws.Range("D41").Value = "'0.09025"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").Value = "'1.481"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").Value = "'0.7705"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").Value = "'0.7214"
$ws.Range("D45").Value = "'15.63"
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "'4.175"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "'0.9985"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "'141.84"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "'1.311"
$ws.Range("E50").Value = "  +14.69%  "
$ws.Range("D51").Value = "'0.00000000371"
$ws.Range("E51").Value = "  +3.04%  "
